# Natmi LR-pairs (Bmp6-Bmpr1b) refresh following Dr Hou advice:
# adds a new "M2" target cluster and recomputes all ligand/receptor
# expression + specificity metrics for the ECs/FAPs/sCs x FAPs/M2/sCs grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    ,@("ECs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 23.061728, 69.18518399999999, 0.6130144106248721, 0.6902769593117909, 3, 1, 1.467406, 4.402218, 0.8178057326420413, 0.8654277364307374, 33.840918037568, 304.5682623381119, 0.5013266992012027, 0.5973848264074955)
    ,@("ECs", "Bmp6", "Bmpr1b", "M2", 3, 1, 23.061728, 69.18518399999999, 0.6130144106248721, 0.6902769593117909, 1, 0.3333333333333333, 0.030706, 0.09211800000000001, 0.01711288002536893, 0.01810938763698815, 0.7081334199679999, 6.373200779712, 0.01049044206284568, 0.01250049303305872)
    ,@("ECs", "Bmp6", "Bmpr1b", "sCs", 3, 1, 23.061728, 69.18518399999999, 0.6130144106248721, 0.6902769593117909, 2, 1, 0.296209, 0.592418, 0.1650813873325899, 0.1164628759322743, 6.831091389151999, 40.98654833491199, 0.1011972693608238, 0.08039163987123665)
    ,@("FAPs", "Bmp6", "Bmpr1b", "FAPs", 3, 1, 1.926013333333334, 5.778040000000001, 0.05119624723650278, 0.05764887294340218, 3, 1, 1.467406, 4.402218, 0.8178057326420413, 0.8654277364307374, 2.826243521413333, 25.43619169272, 0.04186858447977123, 0.04989093361919173)
    ,@("FAPs", "Bmp6", "Bmpr1b", "M2", 3, 1, 1.926013333333334, 5.778040000000001, 0.05119624723650278, 0.05764887294340218, 1, 0.3333333333333333, 0.030706, 0.09211800000000001, 0.01711288002536893, 0.01810938763698815, 0.05914016541333334, 0.5322614887200001, 0.0008761152367073978, 0.001043985786967548)
    ,@("FAPs", "Bmp6", "Bmpr1b", "sCs", 3, 1, 1.926013333333334, 5.778040000000001, 0.05119624723650278, 0.05764887294340218, 2, 1, 0.296209, 0.592418, 0.1650813873325899, 0.1164628759322743, 0.5705024834533334, 3.423014900720001, 0.008451547520024151, 0.006713953537242892)
    ,@("sCs", "Bmp6", "Bmpr1b", "FAPs", 2, 1, 12.632464, 25.264928, 0.3357893421386252, 0.2520741677448068, 3, 1, 1.467406, 4.402218, 0.8178057326420413, 0.8654277364307374, 18.536953468384, 111.221720810304, 0.2746104489610675, 0.2181519764040502)
    ,@("sCs", "Bmp6", "Bmpr1b", "M2", 2, 1, 12.632464, 25.264928, 0.3357893421386252, 0.2520741677448068, 1, 0.3333333333333333, 0.030706, 0.09211800000000001, 0.01711288002536893, 0.01810938763698815, 0.387892439584, 2.327354637504, 0.005746322725815854, 0.004564908816961881)
    ,@("sCs", "Bmp6", "Bmpr1b", "sCs", 2, 1, 12.632464, 25.264928, 0.3357893421386252, 0.2520741677448068, 2, 1, 0.296209, 0.592418, 0.1650813873325899, 0.1164628759322743, 3.741849528976, 14.967398115904, 0.05543257045174195, 0.02935728252379474)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $row = $rowData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

